$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 3772
$ws.Range("E2").Value = 156
$ws.Range("F2").Value = 187
$ws.Range("G2").Value = 93
$ws.Range("H2").Value = 66
$ws.Range("I2").Value = 66
$ws.Range("K2").Value = 4244
$ws.Range("L2").Value = 2953
$ws.Range("M2").Value = 1291
$ws.Range("N2").Value = 1291
$ws.Range("P2").Value = 287
$ws.Range("Q2").Value = 539
$ws.Range("R2").Value = -343
$ws.Range("S2").Value = -237
$ws.Range("T2").Value = 559
$ws.Range("U2").Value = -20
$ws.Range("V2").Value = 2474
$ws.Range("W2").Value = 4.14
$ws.Range("X2").Value = 1.74
$ws.Range("Y2").Value = 4.92
$ws.Range("Z2").Value = 1.51
$ws.Range("AA2").Value = 228.72
$ws.Range("AB2").Value = 354.77
$ws.Range("AC2").Value = 115
$ws.Range("AD2").Value = 26.29
$ws.Range("AE2").Value = 2248
$ws.Range("AF2").Value = 1.34
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 57385305
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# --- Row 3 ---
$ws.Range("D3").Value = 4339
$ws.Range("E3").Value = 240
$ws.Range("F3").Value = 240
$ws.Range("G3").Value = 219
$ws.Range("H3").Value = 198
$ws.Range("I3").Value = 198
$ws.Range("K3").Value = 4291
$ws.Range("L3").Value = 2649
$ws.Range("M3").Value = 1642
$ws.Range("N3").Value = 1642
$ws.Range("P3").Value = 328
$ws.Range("Q3").Value = 652
$ws.Range("R3").Value = -161
$ws.Range("S3").Value = -419
$ws.Range("T3").Value = 485
$ws.Range("U3").Value = 166
$ws.Range("V3").Value = 2031
$ws.Range("W3").Value = 5.54
$ws.Range("X3").Value = 4.56
$ws.Range("Y3").Value = 13.49
$ws.Range("Z3").Value = 4.63
$ws.Range("AA3").Value = 161.31
$ws.Range("AB3").Value = 405.01
$ws.Range("AC3").Value = 313
$ws.Range("AD3").Value = 22.56
$ws.Range("AE3").Value = 2504
$ws.Range("AF3").Value = 2.82
$ws.Range("AG3").Value = 50
$ws.Range("AH3").Value = 0.71
$ws.Range("AI3").Value = 16.58
$ws.Range("AJ3").Value = 65553378
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# --- Row 4 ---
$ws.Range("D4").Value = 4286
$ws.Range("E4").Value = 245
$ws.Range("F4").Value = 245
$ws.Range("G4").Value = 268
$ws.Range("H4").Value = 223
$ws.Range("I4").Value = 223
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 4394
$ws.Range("L4").Value = 2382
$ws.Range("M4").Value = 2013
$ws.Range("N4").Value = 2007
$ws.Range("O4").Value = 6
$ws.Range("P4").Value = 373
$ws.Range("Q4").Value = 605
$ws.Range("R4").Value = -474
$ws.Range("S4").Value = -162
$ws.Range("T4").Value = 493
$ws.Range("U4").Value = 111
$ws.Range("V4").Value = 1735
$ws.Range("W4").Value = 5.71
$ws.Range("X4").Value = 5.19
$ws.Range("Y4").Value = 12.21
$ws.Range("Z4").Value = 5.13
$ws.Range("AA4").Value = 118.34
$ws.Range("AB4").Value = 438.98
$ws.Range("AC4").Value = 320
$ws.Range("AD4").Value = 16.72
$ws.Range("AE4").Value = 2691
$ws.Range("AF4").Value = 1.99
$ws.Range("AG4").Value = 50
$ws.Range("AH4").Value = 0.93
$ws.Range("AI4").Value = 16.74
$ws.Range("AJ4").Value = 74524102

# --- Row 5 ---
$ws.Range("D5").Value = 4617
$ws.Range("E5").Value = 208
$ws.Range("F5").Value = 208
$ws.Range("G5").Value = 182
$ws.Range("H5").Value = 129
$ws.Range("I5").Value = 128
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 4776
$ws.Range("L5").Value = 2701
$ws.Range("M5").Value = 2075
$ws.Range("N5").Value = 2069
$ws.Range("O5").Value = 6
$ws.Range("P5").Value = 373
$ws.Range("Q5").Value = 1092
$ws.Range("R5").Value = -1426
$ws.Range("S5").Value = 330
$ws.Range("T5").Value = 1120
$ws.Range("U5").Value = -28
$ws.Range("V5").Value = 2089
$ws.Range("W5").Value = 4.5
$ws.Range("X5").Value = 2.8
$ws.Range("Y5").Value = 6.3
$ws.Range("Z5").Value = 2.82
$ws.Range("AA5").Value = 130.16
$ws.Range("AB5").Value = 468.25
$ws.Range("AC5").Value = 172
$ws.Range("AD5").Value = 19.99
$ws.Range("AE5").Value = 2774
$ws.Range("AF5").Value = 1.24
$ws.Range("AG5").Value = 50
$ws.Range("AH5").Value = 1.45
$ws.Range("AI5").Value = 29.06
$ws.Range("AJ5").Value = 74524102

# --- Row 6 ---
$ws.Range("D6").Value = 5068
$ws.Range("E6").Value = 121
$ws.Range("F6").Value = 121
$ws.Range("G6").Value = 129
$ws.Range("H6").Value = 33
$ws.Range("I6").Value = 19
$ws.Range("K6").Value = 7443
$ws.Range("L6").Value = 4636
$ws.Range("M6").Value = 2807
$ws.Range("N6").Value = 2045
$ws.Range("P6").Value = 373
$ws.Range("Q6").Value = 736
$ws.Range("R6").Value = -863
$ws.Range("S6").Value = 121
$ws.Range("T6").Value = 633
$ws.Range("U6").Value = 103
$ws.Range("V6").Value = 3786
$ws.Range("W6").Value = 2.38
$ws.Range("X6").Value = 0.65
$ws.Range("Y6").Value = 0.9399999999999999
$ws.Range("Z6").Value = 0.53
$ws.Range("AA6").Value = 165.13
$ws.Range("AB6").Value = 449.63
$ws.Range("AC6").Value = 26
$ws.Range("AD6").Value = 107.62
$ws.Range("AE6").Value = 2743
$ws.Range("AF6").Value = 1.02
$ws.Range("AG6").Value = 25
$ws.Range("AH6").Value = 0.89
$ws.Range("AI6").Value = 96.27
$ws.Range("AJ6").Value = 74524102

# --- Row 7 ---
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# --- Row 8 ---
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# --- Row 9 ---
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
